$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status "En proceso" -> "Cerrada" for the first record (F4)
$ws.Range("F4").Value = "Cerrada"

# Set real closing date (FECHA REAL DE CIERRE) for E4
$ws.Range("E4").Value = (Get-Date -Year 2016 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0).Date

# Update the active selection to F5
$ws.Range("F5").Select()
